$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.2573958041829485
$ws.Range("C2").Value = 0.04748089165519787
$ws.Range("D2").Value = 0.03165052555380043
$ws.Range("E2").Value = 0.1647267180267775
$ws.Range("F2").Value = 0.7170521533020491
$ws.Range("I2").Value = 0.5570539027569694
$ws.Range("K2").Value = 0.2750332023945532
$ws.Range("M2").Value = 0.2202528939603283
$ws.Range("O2").Value = 2.47108682482488

$ws.Range("B3").Value = 0.2256133626909502
$ws.Range("C3").Value = 0.04230195993055474
$ws.Range("D3").Value = 0.02950460241036978
$ws.Range("E3").Value = 0.1539407558026511
$ws.Range("F3").Value = 0.7163095727575524
$ws.Range("I3").Value = 0.5619326161867093
$ws.Range("K3").Value = 0.2399888997756108
$ws.Range("M3").Value = 0.1977350549761994
$ws.Range("O3").Value = 2.48182013221161

$ws.Range("B4").Value = 0.2060647825605599
$ws.Range("C4").Value = 0.03910215581566945
$ws.Range("D4").Value = 0.02817525147519007
$ws.Range("E4").Value = 0.1474358730611982
$ws.Range("F4").Value = 0.7163156781840101
$ws.Range("I4").Value = 0.5652477437201071
$ws.Range("K4").Value = 0.2184091649826883
$ws.Range("M4").Value = 0.1839742000649309
$ws.Range("O4").Value = 2.489921163682126

$ws.Range("B5").Value = 0.1980904729067561
$ws.Range("C5").Value = 0.03779324701719133
$ws.Range("D5").Value = 0.02763060759723857
$ws.Range("E5").Value = 0.1448144618964307
$ws.Range("F5").Value = 0.7164343057102371
$ws.Range("I5").Value = 0.5666789925480202
$ws.Range("K5").Value = 0.209600014440511
$ws.Range("M5").Value = 0.1783829147458462
$ws.Range("O5").Value = 2.493602034970607

$ws.Range("B6").Value = 0.1967658687569269
$ws.Range("C6").Value = 0.03757560534639026
$ws.Range("D6").Value = 0.02753999429141629
$ws.Range("E6").Value = 0.1443809465395418
$ws.Range("F6").Value = 0.7164610168747743
$ws.Range("I6").Value = 0.5669214991463178
$ws.Range("K6").Value = 0.2081363550997395
$ws.Range("M6").Value = 0.1774554746049688
$ws.Range("O6").Value = 2.494236162507875

$ws.Range("B7").Value = 0.205957270444344
$ws.Range("C7").Value = 0.0390845234457089
$ws.Range("D7").Value = 0.02816791801008378
$ws.Range("E7").Value = 0.1474004011274275
$ws.Range("F7").Value = 0.7163168078560176
$ws.Range("I7").Value = 0.5652667209518363
$ws.Range("K7").Value = 0.2182904227716733
$ws.Range("M7").Value = 0.183898727756663
$ws.Range("O7").Value = 2.489969268410391

$ws.Range("B8").Value = 0.2464445435215623
$ws.Range("C8").Value = 0.04569935633921318
$ws.Range("D8").Value = 0.03091306692999041
$ws.Range("E8").Value = 0.1609831416518759
$ws.Range("F8").Value = 0.7167001816052831
$ws.Range("I8").Value = 0.5586697228257513
$ws.Range("K8").Value = 0.2629631383053095
$ws.Range("M8").Value = 0.2124751952889241
$ws.Range("O8").Value = 2.47447396015481

$ws.Range("B9").Value = 0.3255538798997577
$ws.Range("C9").Value = 0.0585115465522108
$ws.Range("D9").Value = 0.03620197935646985
$ws.Range("E9").Value = 0.1885647269041755
$ws.Range("F9").Value = 0.7211214530595029
$ws.Range("I9").Value = 0.5482712584070732
$ws.Range("K9").Value = 0.3500553480227211
$ws.Range("M9").Value = 0.2690345788772603
$ws.Range("O9").Value = 2.456088234992507

$ws.Range("B10").Value = 0.3834852254437351
$ws.Range("C10").Value = 0.06782640898630632
$ws.Range("D10").Value = 0.04002908673394501
$ws.Range("E10").Value = 0.209423772210684
$ws.Range("F10").Value = 0.7266129998833364
$ws.Range("I10").Value = 0.5421825645107212
$ws.Range("K10").Value = 0.4137155257720053
$ws.Range("M10").Value = 0.3109158923864612
$ws.Range("O10").Value = 2.449917869584255

$ws.Range("B11").Value = 0.4097953787275515
$ws.Range("C11").Value = 0.07204248513632194
$ws.Range("D11").Value = 0.04175717666011991
$ws.Range("E11").Value = 0.2190463014861095
$ws.Range("F11").Value = 0.7295997812107515
$ws.Range("I11").Value = 0.5397502699874579
$ws.Range("K11").Value = 0.4426025116471521
$ws.Range("M11").Value = 0.3300421186969089
$ws.Range("O11").Value = 2.44870900046152

$ws.Range("B12").Value = 0.4197517723564204
$ws.Range("C12").Value = 0.07363590661043418
$ws.Range("D12").Value = 0.04240968049538907
$ws.Range("E12").Value = 0.2227096001106261
$ws.Range("F12").Value = 0.7308011475503307
$ws.Range("I12").Value = 0.5388778192227974
$ws.Range("K12").Value = 0.4535304850826378
$ws.Range("M12").Value = 0.337295482299254
$ws.Range("O12").Value = 2.448481378540038

$ws.Range("B13").Value = 0.4176077902124007
$ws.Range("C13").Value = 0.07329287415457486
$ws.Range("D13").Value = 0.04226923655649273
$ws.Range("E13").Value = 0.221919774128267
$ws.Range("F13").Value = 0.7305392824520993
$ws.Range("I13").Value = 0.5390635541027606
$ws.Range("K13").Value = 0.4511774432108382
$ws.Range("M13").Value = 0.335732865720793
$ws.Range("O13").Value = 2.448520159407053

$ws.Range("B14").Value = 0.4106146341191845
$ws.Range("C14").Value = 0.07217363960657508
$ws.Range("D14").Value = 0.04181089654310455
$ws.Range("E14").Value = 0.219347291987475
$ws.Range("F14").Value = 0.7296972082863675
$ws.Range("I14").Value = 0.5396775183547504
$ws.Range("K14").Value = 0.4435017839175828
$ws.Range("M14").Value = 0.3306386432730761
$ws.Range("O14").Value = 2.448685659302356

$ws.Range("B15").Value = 0.406330238239093
$ws.Range("C15").Value = 0.07148766892129288
$ws.Range("D15").Value = 0.04152990338582185
$ws.Range("E15").Value = 0.2177741131480886
$ws.Range("F15").Value = 0.7291905755645089
$ws.Range("I15").Value = 0.5400599214475541
$ws.Range("K15").Value = 0.438798786581998
$ws.Range("M15").Value = 0.3275196755406569
$ws.Range("O15").Value = 2.4488170155646

$ws.Range("B16").Value = 0.3817648855555831
$ws.Range("C16").Value = 0.06755044540850008
$ws.Range("D16").Value = 0.03991588988633765
$ws.Range("E16").Value = 0.2087976268839071
$ws.Range("F16").Value = 0.7264276455719454
$ws.Range("I16").Value = 0.542348318261979
$ws.Range("K16").Value = 0.4118261948931661
$ws.Range("M16").Value = 0.3096674422559502
$ws.Range("O16").Value = 2.450029059170419

$ws.Range("B17").Value = 0.3666834462306952
$ws.Range("C17").Value = 0.06512959343974956
$ws.Range("D17").Value = 0.03892242157712644
$ws.Range("E17").Value = 0.2033252420598117
$ws.Range("F17").Value = 0.7248578823874539
$ws.Range("I17").Value = 0.5438386677224756
$ws.Range("K17").Value = 0.395260498691556
$ws.Range("M17").Value = 0.2987346857489257
$ws.Range("O17").Value = 2.451182160968727

$ws.Range("B18").Value = 0.3580049731122585
$ws.Range("C18").Value = 0.0637351815080649
$ws.Range("D18").Value = 0.03834979454724419
$ws.Range("E18").Value = 0.2001902361395054
$ws.Range("F18").Value = 0.7240009859789964
$ws.Range("I18").Value = 0.5447276382183261
$ws.Range("K18").Value = 0.3857255682303276
$ws.Range("M18").Value = 0.2924534377089856
$ws.Range("O18").Value = 2.451995790778938

$ws.Range("B19").Value = 0.3550659144170822
$ws.Range("C19").Value = 0.06326271534358341
$ws.Range("D19").Value = 0.03815570604052709
$ws.Range("E19").Value = 0.1991309278267153
$ws.Range("F19").Value = 0.7237187522338289
$ws.Range("I19").Value = 0.5450340809793843
$ws.Range("K19").Value = 0.3824960553072003
$ws.Range("M19").Value = 0.2903279144329716
$ws.Range("O19").Value = 2.452297091812142

$ws.Range("B20").Value = 0.3682893116560138
$ws.Range("C20").Value = 0.06538750483291267
$ws.Range("D20").Value = 0.03902830348705777
$ws.Range("E20").Value = 0.2039064846068115
$ws.Range("F20").Value = 0.7250202260695957
$ws.Range("I20").Value = 0.5436767297858829
$ws.Range("K20").Value = 0.3970246506636101
$ws.Range("M20").Value = 0.2998977734480164
$ws.Range("O20").Value = 2.451043843320463

$ws.Range("B21").Value = 0.4126688790663593
$ws.Range("C21").Value = 0.07250247087941375
$ws.Range("D21").Value = 0.04194557358007245
$ws.Range("E21").Value = 0.2201023624359166
$ws.Range("F21").Value = 0.729942636401411
$ws.Range("I21").Value = 0.5394958623971284
$ws.Range("K21").Value = 0.4457566109228139
$ws.Range("M21").Value = 0.3321346500795528
$ws.Range("O21").Value = 2.44863079895876

$ws.Range("B22").Value = 0.4416342682096968
$ws.Range("C22").Value = 0.07713433162842875
$ws.Range("D22").Value = 0.04384116759127465
$ws.Range("E22").Value = 0.2308008174481984
$ws.Range("F22").Value = 0.7335697260520391
$ws.Range("I22").Value = 0.537046774560082
$ws.Range("K22").Value = 0.4775419941501866
$ws.Range("M22").Value = 0.35326554775191
$ws.Range("O22").Value = 2.44839539055522

$ws.Range("B23").Value = 0.4261786458274344
$ws.Range("C23").Value = 0.07466390100887565
$ws.Range("D23").Value = 0.04283047211076507
$ws.Range("E23").Value = 0.2250803857137242
$ws.Range("F23").Value = 0.7315963412948605
$ws.Range("I23").Value = 0.5383279469126201
$ws.Range("K23").Value = 0.4605835357411081
$ws.Range("M23").Value = 0.3419818934948538
$ws.Range("O23").Value = 2.448398155226812

$ws.Range("B24").Value = 0.3675633246329539
$ws.Range("C24").Value = 0.06527091128744189
$ws.Range("D24").Value = 0.03898043884796465
$ws.Range("E24").Value = 0.2036436701688231
$ws.Range("F24").Value = 0.724946688512297
$ws.Range("I24").Value = 0.543749841789996
$ws.Range("K24").Value = 0.3962271120681464
$ws.Range("M24").Value = 0.2993719286098013
$ws.Range("O24").Value = 2.451105907339979

$ws.Range("B25").Value = 0.3041849854037935
$ws.Range("C25").Value = 0.05506269063798186
$ws.Range("D25").Value = 0.0347814059239866
$ws.Range("E25").Value = 0.181000031485965
$ws.Range("F25").Value = 0.7195318920116591
$ws.Range("I25").Value = 0.5508122023643551
$ws.Range("K25").Value = 0.3265508511610449
$ws.Range("M25").Value = 0.2536770058463276
$ws.Range("O25").Value = 2.459775065628975
